$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(44186, "Bánovce nad Bebravou", 15),
  @(44186, "Banská Bystrica", 57),
  @(44186, "Banská Štiavnica", 2),
  @(44186, "Bardejov", 14),
  @(44186, "Bratislava", 296),
  @(44186, "Brezno", 14),
  @(44186, "Bytča", 21),
  @(44186, "Čadca", 63),
  @(44186, "Detva", 3),
  @(44186, "Dolný Kubín", 17),
  @(44186, "Dunajská Streda", 15),
  @(44186, "Galanta", 20),
  @(44186, "Gelnica", 14),
  @(44186, "Hlohovec", 33),
  @(44186, "Humenné", 31),
  @(44186, "Ilava", 24),
  @(44186, "Kežmarok", 11),
  @(44186, "Komárno", 8),
  @(44186, "Košice", 133),
  @(44186, "Košice - okolie", 74),
  @(44186, "Krupina", 5),
  @(44186, "Kysucké Nové Mesto", 41),
  @(44186, "Levice", 7),
  @(44186, "Levoča", 13),
  @(44186, "Liptovský Mikuláš", 19),
  @(44186, "Lučenec", 63),
  @(44186, "Malacky", 27),
  @(44186, "Martin", 106),
  @(44186, "Medzilaborce", 5),
  @(44186, "Michalovce", 32),
  @(44186, "Myjava", 40),
  @(44186, "Námestovo", 14),
  @(44186, "Nitra", 32),
  @(44186, "Nové Mesto nad Váhom", 28),
  @(44186, "Nové Zámky", 47),
  @(44186, "Partizánske", 7),
  @(44186, "Pezinok", 17),
  @(44186, "Piešťany", 48),
  @(44186, "Poltár", 4),
  @(44186, "Poprad", 25),
  @(44186, "Považská Bystrica", 12),
  @(44186, "Prešov", 165),
  @(44186, "Prievidza", 44),
  @(44186, "Púchov", 4),
  @(44186, "Revúca", 49),
  @(44186, "Rimavská Sobota", 27),
  @(44186, "Rožňava", 5),
  @(44186, "Ružomberok", 34),
  @(44186, "Sabinov", 31),
  @(44186, "Senec", 52),
  @(44186, "Senica", 52),
  @(44186, "Skalica", 40),
  @(44186, "Snina", 7),
  @(44186, "Sobrance", 8),
  @(44186, "Spišská Nová Ves", 66),
  @(44186, "Stará Ľubovňa", 39),
  @(44186, "Stropkov", 15),
  @(44186, "Svidník", 37),
  @(44186, "Šaľa", 9),
  @(44186, "Topoľčany", 17),
  @(44186, "Trebišov", 73),
  @(44186, "Trenčín", 74),
  @(44186, "Trnava", 120),
  @(44186, "Turčianske Teplice", 7),
  @(44186, "Tvrdošín", 20),
  @(44186, "Veľký Krtíš", 8),
  @(44186, "Vranov nad Topľou", 39),
  @(44186, "Zlaté Moravce", 8),
  @(44186, "Zvolen", 44),
  @(44186, "Žarnovica", 1),
  @(44186, "Žiar nad Hronom", 9),
  @(44186, "Žilina", 102)
)

$startRow = 6567
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
  $ws.Cells.Item($r, 3).Value = $data[$i][2]
}


